$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert 6 new columns before old column L (fix_time),
# shifting fix_time..target_time (and everything after) right by 6.
$ws.Range("L1:Q1").EntireColumn.Insert()

# Step 2: insert 3 new columns before the old column R (prime_traj),
# which is now at X after step 1. This shifts prime_traj.. right by 3 more.
$ws.Range("X1:Z1").EntireColumn.Insert()

# --- Row 1 (header names) for newly inserted columns ---
$ws.Range("L1").Value = "fix_duration"
$ws.Range("M1").Value = "mask1_duration"
$ws.Range("N1").Value = "mask2_duration"
$ws.Range("O1").Value = "prime_duration"
$ws.Range("P1").Value = "mask3_duration"
$ws.Range("Q1").Value = "target_duration"
$ws.Range("X1").Value = "categor_time"
$ws.Range("Y1").Value = "recog_time"
$ws.Range("Z1").Value = "pas_time"

# --- Row 2 (header descriptions) for newly inserted columns ---
$ws.Range("L2").Value = "Fixation duration"
$ws.Range("M2").Value = "Long mask duration"
$ws.Range("N2").Value = "Forward mask duration"
$ws.Range("O2").Value = "Prime duration"
$ws.Range("P2").Value = "Backward mask duration"
$ws.Range("Q2").Value = "target duration"
$ws.Range("X2").Value = "timestamp when categorization task was displayed"
$ws.Range("Y2").Value = "timestamp when recognition task was displayed"
$ws.Range("Z2").Value = "timestamp when PAS task was displayed"

# --- Row 2 descriptions for the shifted fix_time..target_time columns ---
# (these now describe timestamps rather than durations)
$ws.Range("R2").Value = "timestamp when fixation was displayed"
$ws.Range("S2").Value = "timestamp when mask1 was displayed"
$ws.Range("T2").Value = "timestamp when mask2 was displayed"
$ws.Range("U2").Value = "timestamp when prime was displayed"
$ws.Range("V2").Value = "timestamp when mask3 was displayed"
$ws.Range("W2").Value = "timestamp when target was displayed"

# Restore the active cell selection to match the edited workbook state
$ws.Range("Z3").Select() | Out-Null
